# Sync Excel worksheet content/view state to match the authoritative SVN copy.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The code-list cell used literal numeric codes ("30000,30001"); replace it with
# the human-readable Chinese layer labels ("地表层;地板层").
$ws.Range("J3").Value = "地表层;地板层"

# Restore the saved cursor position: the active cell/selection moves from L9 to J6.
[void]$ws.Range("J6").Select()
